$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Range("A1").Value = "file_name"
$ws.Range("B1").Value = "send_email_to"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Completion date"

# ---- Data rows ----
$ws.Range("A2").Value = "Mihai Popescu"
$ws.Range("B2").Value = "mihai.popescu@gmail.com"
$ws.Range("C2").Value = "Mihai Popescu"
$ws.Range("D2").Value = "December 14th 2022"

$ws.Range("A3").Value = "Dragos Ionescu"
$ws.Range("B3").Value = "dragos.ionescu@gmail.com"
$ws.Range("C3").Value = "Dragos Ionescu"
$ws.Range("D3").Value = "December 14th 2022"

$ws.Range("A4").Value = "Alex Georgescu"
$ws.Range("B4").Value = "alex.georgescu@gmail.com"
$ws.Range("C4").Value = "Alex Georgescu"
$ws.Range("D4").Value = "December 14th 2022"
